# JetEstate.xlsx edit
# - Make "Specs" the active sheet/tab (was "Preview - Index")
# - Select cell E10 on "Specs"
# - Re-order a couple of "Pages" entries on "Specs" (About / Properties / Search
#   Result / Item Page) so that "Properties" (and its "Search Result" child)
#   now sits above "About" (and its "Item Page" child)
# - Keep the existing B21:L21 selection on "Preview - Index", it just stops
#   being the active/visible tab

$wb = $excel.ActiveWorkbook

$specs = $wb.Worksheets.Item("Specs")
$preview = $wb.Worksheets.Item("Preview - Index")

# --- Re-shuffle the four "pages" rows on the Specs sheet ---------------
$about        = $specs.Range("A7").Value()
$properties   = $specs.Range("A8").Value()
$searchResult = $specs.Range("B9").Value()
$itemPage     = $specs.Range("B10").Value()

$specs.Range("A7").Value = $properties
$specs.Range("B8").Value = $searchResult
$specs.Range("A8").ClearContents()
$specs.Range("B9").Value = $itemPage
$specs.Range("A10").Value = $about
$specs.Range("B10").ClearContents()

# --- Switch the active tab from "Preview - Index" to "Specs" -----------
[void]$specs.Activate()
[void]$specs.Range("E10").Select()

# "Preview - Index" keeps its previous selection, it's just no longer the
# active sheet.
[void]$preview.Range("B21:L21").Select()
[void]$specs.Activate()
